$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" value from EA-8 to EA-7
$ws.Range("B9").Value = "EA-7"
$ws.Range("C9").Value = "EA-7"

# Delete row 25 (LOB1240 requirement) entirely, shifting cells up
$ws.Rows(25).Delete()
